# Apply weekly update: insert 3 new price rows at the top of the
# "Mandarina" data block (row 676), shifting the existing rows down by 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 676 (this pushes old rows 676-743 to 679-746)
$ws.Rows("676:678").Insert()

# Common (constant) column values reused across this data block
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$prodId    = 100102
$prod      = "Cítricos"
$catId     = 100102004
$cat       = "Mandarina"
$unidad    = "`$/bandeja 10 kilos"
$origen    = "Provincia de Limarí"
$kgUnidad  = 10

function Set-PrecioRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg) {
    $ws.Cells.Item($Row, 1).Value  = $mercadoId
    $ws.Cells.Item($Row, 2).Value  = $mercado
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $prodId
    $ws.Cells.Item($Row, 8).Value  = $prod
    $ws.Cells.Item($Row, 9).Value  = $catId
    $ws.Cells.Item($Row, 10).Value = $cat
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $unidad
    $ws.Cells.Item($Row, 18).Value = $origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $kgUnidad
}

# New rows for the latest weekly report (Fecha serial 44858 = 2022-10-24)
Set-PrecioRow 676 44858 "Murcott" "Especial" 500 6000 6500 6250 625
Set-PrecioRow 677 44858 "Murcott" "Primera"  508 5000 5500 5254 525
Set-PrecioRow 678 44858 "Murcott" "Segunda"  360 4000 4500 4250 425
